$wb = $excel.ActiveWorkbook

# Sheet references
$wsComparativa = $wb.Worksheets.Item("comparativa de las 3")

# --- sheet4 ("comparativa de las 3"): remove the old row-index column (A) ---
$wsComparativa.Range("A7:A29").ClearContents()

# --- sheet4: add the new title banner in B1:H3 (merged) ---
$titleCell = $wsComparativa.Range("B1")
$titleCell.Value = "Comparativa dentro de la misma gráfica para cada uno de los parámetros estudiados"
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 20
$titleCell.HorizontalAlignment = -4108
$titleCell.VerticalAlignment = -4108
$wsComparativa.Range("B1:H3").Merge() | Out-Null

# --- sheet4: page setup (paper size / orientation) ---
$wsComparativa.PageSetup.PaperSize = 9
$wsComparativa.PageSetup.Orientation = 1

# --- chart title fix: the 4th chart ("Gráfico 4") on sheet4 actually plots the
#     "Eficiencia" data block (rows 88:111) but was mislabeled "Coste" ---
foreach ($co in $wsComparativa.ChartObjects()) {
    if ($co.Name -eq "Gráfico 4") {
        $co.Chart.ChartTitle.Text = "Eficiencia (comparativa de 3 casos)"
    }
}

# --- make sheet4 the active tab, with the new selection on the cleared range ---
$wsComparativa.Activate()
$wsComparativa.Range("A7:A29").Select()
